# Update price list workbook: refresh date stamp and price values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Date stamp in A1 (merged A1:D1) - moved one month forward (serial date value)
$ws.Range("A1").Value = 45436

# SOPORTE VISILLO c/Tornillo prices
$ws.Range("D23").Value = 48.442
$ws.Range("D24").Value = 48.442

# SOPORTE VISILLO a Presion prices
$ws.Range("D25").Value = 28.707
$ws.Range("D26").Value = 28.707

# Soporte de FLEJE prices
$ws.Range("D41").Value = 106.967
$ws.Range("D42").Value = 106.967
$ws.Range("D43").Value = 121.996
$ws.Range("D44").Value = 152.625
